$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.833.18"
$ws.Range("E2").Value = "  -0.77%  "

# Row 3
$ws.Range("D3").Value = "2.918.58"
$ws.Range("E3").Value = "  +0.21%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "356.32"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.39"
$ws.Range("E6").Value = "  -2.71%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.566"
$ws.Range("E7").Value = "  +1.34%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.623"
$ws.Range("E9").Value = "  -1.47%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.11"
$ws.Range("E10").Value = "  -2.13%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0871"
$ws.Range("E11").Value = "  +1.20%  "

# Row 12
$ws.Range("E12").Value = "  +1.07%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.52"
$ws.Range("E13").Value = "  -2.53%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.81"
$ws.Range("E14").Value = "  -0.04%  "

# Row 15
$ws.Range("D15").Value = "3.378.48"
$ws.Range("E15").Value = "  +0.23%  "

# Row 16
$ws.Range("D16").Value = "2.909.14"
$ws.Range("E16").Value = "  +0.37%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.983"
$ws.Range("E17").Value = "  -1.62%  "

# Row 18
$ws.Range("D18").Value = "51.814.60"
$ws.Range("E18").Value = "  -0.86%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.37"
$ws.Range("E19").Value = "  +1.83%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.54"
$ws.Range("E20").Value = "  -1.36%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.93"
$ws.Range("E21").Value = "  -1.85%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0979"
$ws.Range("E22").Value = "  -0.25%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.50"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.71"
$ws.Range("E24").Value = "  -0.43%  "

# Row 25
$ws.Range("E25").Value = "  +0.31%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.186"
$ws.Range("E26").Value = "  +12.72%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.72"
$ws.Range("E27").Value = "  +18.50%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "26.90"
$ws.Range("E28").Value = "  +0.54%  "

# Row 29
$ws.Range("E29").Value = "  +0.11%  "

# Row 30
$ws.Range("E30").Value = "  +10.83%  "

# Row 31
$ws.Range("E31").Value = "  -0.81%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.52"
$ws.Range("E32").Value = "  -0.34%  "

# Row 33
$ws.Range("E33").Value = "  -1.35%  "

# Row 34
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "52.17"
$ws.Range("E34").Value = "  -2.14%  "

# Row 35
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.05"
$ws.Range("E35").Value = "  -8.88%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0442"
$ws.Range("E36").Value = "  -2.44%  "

# Row 37
$ws.Range("E37").Value = "  -0.02%  "

# Row 38
$ws.Range("E38").Value = "  -2.81%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.25"
$ws.Range("E39").Value = "  -2.55%  "

# Row 40
$ws.Range("E40").Value = "  -3.54%  "

# Row 41
$ws.Range("E41").Value = "  -3.15%  "

# Row 42
$ws.Range("E42").Value = "  +2.09%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.76"
$ws.Range("E43").Value = "  -3.39%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.27"
$ws.Range("E44").Value = "  -1.79%  "

# Row 45
$ws.Range("E45").Value = "  +0.13%  "

# Row 46
$ws.Range("E46").Value = "  -1.45%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.47"
$ws.Range("E47").Value = "  -5.79%  "

# Row 48
$ws.Range("D48").Value = "2.123.84"
$ws.Range("E48").Value = "  -3.47%  "

# Row 49
$ws.Range("E49").Value = "  -4.68%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0338"
$ws.Range("E50").Value = "  +0.36%  "

# Row 51
$ws.Range("B51").Value = "SEI"
$ws.Range("C51").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.911"
$ws.Range("E51").Value = "  -5.44%  "
